$wb = $excel.ActiveWorkbook

$wsID4 = $wb.Worksheets.Item("ID 4")
$wsID5 = $wb.Worksheets.Item("ID 5")
$wsID6 = $wb.Worksheets.Item("ID 6")
$wsID7 = $wb.Worksheets.Item("ID 7")
$wsDunk = $wb.Worksheets.Item("Dunk")

# Remove plotted data from "ID 7" sheet (Feed_SiO2 C, Feed_Cl_hach D, Perm_SiO2 I, Perm_Cl_hach J)
$wsID7.Range("C2:D6").ClearContents()
$wsID7.Range("I3:J5").ClearContents()

# Update selections to reflect what was left active after the edit
$wsID4.Range("C30").Select()
$wsID7.Range("C2:J7").Select()

# Activate "ID 6" sheet, which becomes the selected/active tab
$wsID6.Activate()

$wb.Save()
